$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.562.33"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "1.902.54"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'239.09"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'0.4731"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("D8").Value = "'0.2855"
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "'0.06640"
$ws.Range("E9").Value = "  -0.77%  "
$ws.Range("D10").Value = "'19.65"
$ws.Range("E10").Value = "  +4.59%  "
$ws.Range("D11").Value = "'99.89"
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("D12").Value = "'0.07805"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "1.909.41"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "'5.185"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "'0.6752"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'283.80"
$ws.Range("E16").Value = "  +8.40%  "
$ws.Range("D17").Value = "30.565.57"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "'1.002"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.000007470"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.155.93"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'12.71"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'5.404"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "'6.270"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "'9.356"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("D26").Value = "'166.96"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").Value = "'19.27"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("D28").Value = "'2.023"
$ws.Range("E28").Value = "  -2.16%  "
$ws.Range("D29").Value = "'1.380"
$ws.Range("D30").Value = "'0.09937"
$ws.Range("E30").Value = "  -1.76%  "
$ws.Range("D31").Value = "'4.505"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "'1.509"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "'4.253"
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "'0.04749"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "'0.7207"
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("D36").Value = "'1.107"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.725"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01897"
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'6.700"
$ws.Range("E39").Value = "  +6.72%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.569"
$ws.Range("E40").Value = "  -0.79%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'73.84"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.986"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'0.8700"
$ws.Range("E43").Value = "  +1.08%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'104.95"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "'0.4265"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "'986.01"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").Value = "'7.378"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.218"
$ws.Range("E49").Value = "  +4.42%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.1184"
$ws.Range("E50").Value = "  -0.63%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'34.44"
$ws.Range("E51").Value = "  -1.48%  "
